# Add the "intervention_type" column (K) to the worksheet, mirroring the
# header style used by the other header cells (A1:J1), and fill in the
# per-row intervention type values for rows 2-20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell K1 --------------------------------------------------
# Copy the formatting of the existing header cell J1 onto K1 so it picks
# up the same bold/centered/bordered header style, then set its text.
$ws.Range("J1").Copy() | Out-Null
$ws.Range("K1").PasteSpecial(-4122) | Out-Null
$ws.Range("K1").Value = "intervention_type"

# --- Data rows 2-20 ----------------------------------------------------
$values = @(
    "PROCEDURE",
    "DRUG",
    "DRUG",
    "DRUG",
    "OTHER",
    "DIAGNOSTIC_TEST",
    "DIAGNOSTIC_TEST",
    "DRUG",
    "OTHER",
    "DEVICE",
    "OTHER",
    "DIETARY_SUPPLEMENT",
    "DRUG",
    "DEVICE",
    "PROCEDURE",
    "DRUG",
    "DIAGNOSTIC_TEST",
    "DRUG",
    "OTHER"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 11).Value = $values[$i]
}
